# "Generate Report for Handback"
# The CI run that produced this workbook has now received handback
# translations for both locales, so the status report is refreshed:
#   - Overview sheet: status text updated for each locale column
#   - zh-cn sheet: handback target/file recorded (handback datetime was
#     already real, only the status text changes)
#   - de-de sheet: handback target/file/datetime recorded

$wb = $excel.ActiveWorkbook

$ovWs   = $wb.Worksheets.Item("Overview")
$zhWs   = $wb.Worksheets.Item("zh-cn")
$deWs   = $wb.Worksheets.Item("de-de")

$mdFileName = "d169085d-ee20-480a-8fed-8c9db05c8fc5.md"
$mdHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8a910c1f8672e4b630ffa55a9803d3f72009f08/e2e/d169085d-ee20-480a-8fed-8c9db05c8fc5.md"

$zhXlfFileName = "d169085d-ee20-480a-8fed-8c9db05c8fc5.2e8e0b8f07559529eb2e026432081d639ed7e6ef.zh-cn.xlf"
$deXlfFileName = "d169085d-ee20-480a-8fed-8c9db05c8fc5.2e8e0b8f07559529eb2e026432081d639ed7e6ef.de-de.xlf"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status cells move from "ready for
# handoff" to "handed back" now that both translations are in.
# ---------------------------------------------------------------------
$ovWs.Range("E2").Value = $newStatus
$ovWs.Range("F2").Value = $newStatus

# Column widths auto-grow to fit the longer status text.
$ovWs.Columns.Item(5).ColumnWidth = 29.144371396019366
$ovWs.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------
# zh-cn sheet: handback file + target are now known; status text
# refreshed as well. The handback datetime cell already held a real
# timestamp and is refreshed to the latest handback run.
# ---------------------------------------------------------------------
$zhWs.Range("C2").Value = $newStatus
$zhWs.Range("I2").Value = $mdFileName
$zhWs.Hyperlinks.Add($zhWs.Range("I2"), $mdHyperlinkUrl, "", "", $mdFileName) | Out-Null
$zhWs.Range("J2").Value = $zhXlfFileName
$zhWs.Range("K2").Value = "2016-08-27 14:57:26"

$zhWs.Columns.Item(3).ColumnWidth = 29.144371396019366
$zhWs.Columns.Item(9).ColumnWidth = 39.16666666666667
$zhWs.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet: handback file + target + handback datetime are now
# known; status text refreshed as well.
# ---------------------------------------------------------------------
$deWs.Range("C2").Value = $newStatus
$deWs.Range("I2").Value = $mdFileName
$deWs.Hyperlinks.Add($deWs.Range("I2"), $mdHyperlinkUrl, "", "", $mdFileName) | Out-Null
$deWs.Range("J2").Value = $deXlfFileName
$deWs.Range("K2").Value = "2016-08-27 14:57:33"

$deWs.Columns.Item(3).ColumnWidth = 29.144371396019366
$deWs.Columns.Item(9).ColumnWidth = 39.16666666666667
$deWs.Columns.Item(10).ColumnWidth = 39.16666666666667
